$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C15").Value = 8547
$ws.Range("C16:C28").Value = 8414
$ws.Range("C29:C52").Value = 7925
$ws.Range("C53:C71").Value = 7900
